# Applies the updated statistical-description values to sheet1.
# The workbook is a single-sheet "statistical description" table where
# rows 2-17 hold per-field stats (Mean, STD, Min, 1st Quartile, Median,
# 3rd Quartile, Max) in columns C-I. This script writes the refreshed
# values coming from the upstream recalculation described in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C2"  = 548.4950573743329
    "D2"  = 134.1831272940331
    "H2"  = 614

    "C3"  = 40.53816140791415
    "D3"  = 4.708294320196429
    "G3"  = 39.86
    "H3"  = 43.08

    "C4"  = 1.390171292267375
    "D4"  = 2.205551596094132
    "G4"  = 1
    "H4"  = 1.78

    "C5"  = 323.5475284446826
    "D5"  = 10.53104554449298
    "F5"  = 317.48
    "G5"  = 325.43
    "H5"  = 332.29

    "C6"  = 20.68353704656674
    "D6"  = 2.410231282285922
    "F6"  = 19.32
    "G6"  = 20.42
    "H6"  = 22.1

    "C7"  = -76.482122217163
    "D7"  = 22.48240524190338

    "C8"  = 7.699387498606519
    "D8"  = 6.894939634500335

    "C9"  = 9.323020592252917
    "D9"  = 1.687916756027434

    "C10" = 867.830306957395
    "D10" = 0.4611338127051571

    "C11" = 0.5569145285120393
    "D11" = 0.5904721468873692

    "C12" = 22.68993522093406
    "D12" = 12.27851378189118

    "C13" = 0.6713636653473758
    "D13" = 0.7481759703058133

    "C14" = 1.826502029229162
    "D14" = 1.666587381203018

    "C15" = 93.88212221716277
    "D15" = 22.48240524190339

    "C16" = -85.76217271890657
    "D16" = 20.23167260634411
    "E16" = -141.2123840191425
    "G16" = -85.5149694202523
    "H16" = -68.16954289279533

    "C17" = -78.06278522030004
    "D17" = 24.89797032201137
    "E17" = -146.5153820769794
    "G17" = -75.3175485570292
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
